$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while preserving its original
# cell style (forcing text interpretation for values that would
# otherwise be auto-converted to numbers, e.g. "536.01").
function Set-TextValue {
    param($cell, $text)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '58.949.43'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '2.502.36'
$ws.Range('E3').Value = '  -0.69%  '
$ws.Range('E4').Value = '  +0.21%  '
Set-TextValue $ws.Range('D5') '536.01'
$ws.Range('E5').Value = '  -0.18%  '
Set-TextValue $ws.Range('D6') '138.24'
$ws.Range('E6').Value = '  -1.02%  '
$ws.Range('E7').Value = '  -0.07%  '
Set-TextValue $ws.Range('D8') '0.567'
$ws.Range('E8').Value = '  +1.02%  '
$ws.Range('D9').Value = '2.527.73'
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('E10').Value = '  +1.17%  '
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('E12').Value = '  -1.45%  '
$ws.Range('E13').Value = '  -2.03%  '
$ws.Range('D14').Value = '2.985.73'
$ws.Range('E14').Value = '  +0.66%  '
Set-TextValue $ws.Range('D15') '23.27'
$ws.Range('E15').Value = '  +1.96%  '
$ws.Range('D16').Value = '58.992.47'
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '2.517.12'
$ws.Range('E18').Value = '  -1.33%  '
Set-TextValue $ws.Range('D19') '11.12'
$ws.Range('E19').Value = '  +1.87%  '
Set-TextValue $ws.Range('D20') '4.27'
$ws.Range('E20').Value = '  +1.22%  '
Set-TextValue $ws.Range('D21') '325.76'
$ws.Range('E21').Value = '  +1.21%  '
Set-TextValue $ws.Range('D22') '1.01'
$ws.Range('E22').Value = '  +0.65%  '
Set-TextValue $ws.Range('D23') '5.90'
$ws.Range('E23').Value = '  +1.33%  '
Set-TextValue $ws.Range('D24') '64.88'
$ws.Range('E24').Value = '  +4.63%  '
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('E27').Value = '  +1.46%  '
Set-TextValue $ws.Range('D28') '7.64'
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('D29').Value = '0.0₃0779'
$ws.Range('E29').Value = '  +1.77%  '
Set-TextValue $ws.Range('D30') '6.73'
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('E31').Value = '  -1.50%  '
Set-TextValue $ws.Range('D32') '167.90'
$ws.Range('E32').Value = '  +4.39%  '
Set-TextValue $ws.Range('D33') '1.19'
$ws.Range('E33').Value = '  +5.79%  '
$ws.Range('E34').Value = '  +0.01%  '
Set-TextValue $ws.Range('D35') '1.41'
$ws.Range('E35').Value = '  -3.34%  '
Set-TextValue $ws.Range('D36') '18.61'
$ws.Range('E36').Value = '  +0.61%  '
Set-TextValue $ws.Range('D37') '4.12'
$ws.Range('E37').Value = '  -1.93%  '
Set-TextValue $ws.Range('D38') '1.57'
$ws.Range('E38').Value = '  -0.58%  '
Set-TextValue $ws.Range('D39') '36.77'
$ws.Range('E39').Value = '  -0.56%  '
Set-TextValue $ws.Range('D40') '0.832'
$ws.Range('E40').Value = '  +3.70%  '
$ws.Range('E41').Value = '  +0.30%  '
Set-TextValue $ws.Range('D42') '5.29'
$ws.Range('E42').Value = '  +1.37%  '
Set-TextValue $ws.Range('D43') '282.69'
$ws.Range('E43').Value = '  -0.17%  '
$ws.Range('E44').Value = '  -0.40%  '
Set-TextValue $ws.Range('D45') '0.607'
$ws.Range('E45').Value = '  +1.84%  '
Set-TextValue $ws.Range('D46') '130.61'
$ws.Range('E46').Value = '  +6.78%  '
$ws.Range('E47').Value = '  +0.26%  '
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('E49').Value = '  +0.46%  '
Set-TextValue $ws.Range('D50') '0.0225'
$ws.Range('E50').Value = '  +0.84%  '
Set-TextValue $ws.Range('D51') '17.44'
$ws.Range('E51').Value = '  +0.18%  '
